$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be forced
# to Text number format before assignment, otherwise Excel auto-converts
# the string into a real number (e.g. "5.930" -> 5.93, "241.68" -> 241.68
# as a float instead of text), which would not match the source data.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D15",
    "D17",
    "D18",
    "D19",
    "D20",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D31",
    "D32",
    "D33",
    "D34",
    "D36",
    "D37",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.477.50'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').Value = '1.874.83'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '0.7152'
$ws.Range('E5').Value = '  +2.53%  '
$ws.Range('D6').Value = '241.68'
$ws.Range('E6').Value = '  +1.95%  '
$ws.Range('D7').Value = '0.9989'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').Value = '0.07952'
$ws.Range('E8').Value = '  +0.94%  '
$ws.Range('D9').Value = '0.3106'
$ws.Range('E9').Value = '  +3.08%  '
$ws.Range('D10').Value = '25.33'
$ws.Range('E10').Value = '  +6.31%  '
$ws.Range('D11').Value = '0.08264'
$ws.Range('E11').Value = '  +1.62%  '
$ws.Range('D12').Value = '0.7305'
$ws.Range('E12').Value = '  +3.79%  '
$ws.Range('D14').Value = '1.862.58'
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('D15').Value = '91.12'
$ws.Range('E15').Value = '  +2.09%  '
$ws.Range('D16').Value = '29.481.05'
$ws.Range('E16').Value = '  +0.89%  '
$ws.Range('D17').Value = '5.930'
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').Value = '245.93'
$ws.Range('E18').Value = '  +4.56%  '
$ws.Range('D19').Value = '0.000007880'
$ws.Range('D20').Value = '13.35'
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('D21').Value = '2.115.91'
$ws.Range('E21').Value = '  -1.03%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').Value = '8.019'
$ws.Range('E22').Value = '  +6.87%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '0.9985'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').Value = '0.9993'
$ws.Range('D25').Value = '0.1613'
$ws.Range('E25').Value = '  +14.27%  '
$ws.Range('D26').Value = '163.30'
$ws.Range('E26').Value = '  +0.66%  '
$ws.Range('D27').Value = '9.042'
$ws.Range('E27').Value = '  +1.90%  '
$ws.Range('D28').Value = '18.31'
$ws.Range('E28').Value = '  +1.78%  '
$ws.Range('D29').Value = '1.355'
$ws.Range('E29').Value = '  -3.09%  '
$ws.Range('E30').Value = '  +1.10%  '
$ws.Range('D31').Value = '4.383'
$ws.Range('E31').Value = '  +2.87%  '
$ws.Range('D32').Value = '4.117'
$ws.Range('E32').Value = '  +2.64%  '
$ws.Range('D33').Value = '0.05268'
$ws.Range('E33').Value = '  +2.32%  '
$ws.Range('D34').Value = '1.951'
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('E35').Value = '  +2.97%  '
$ws.Range('D36').Value = '0.7276'
$ws.Range('E36').Value = '  +3.06%  '
$ws.Range('D37').Value = '2.674'
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('E38').Value = '  +1.50%  '
$ws.Range('D39').Value = '1.217.91'
$ws.Range('E39').Value = '  +6.25%  '
$ws.Range('D40').Value = '2.702'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').Value = '0.9128'
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').Value = '73.89'
$ws.Range('E42').Value = '  +5.69%  '
$ws.Range('D43').Value = '6.113'
$ws.Range('E43').Value = '  +2.98%  '
$ws.Range('D44').Value = '0.9994'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').Value = '102.24'
$ws.Range('E45').Value = '  -0.58%  '
$ws.Range('D46').Value = '0.5288'
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('D47').Value = '2.011.44'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').Value = '1.796'
$ws.Range('E48').Value = '  +3.59%  '
$ws.Range('D49').Value = '2.918'
$ws.Range('E49').Value = '  +9.01%  '
$ws.Range('D50').Value = '9.338'
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('D51').Value = '0.4324'
$ws.Range('E51').Value = '  +2.30%  '
